# A new weekly price record was inserted at row 297 ("Fecha" 45124, i.e. 2023-07-17).
# All the existing records that were in rows 297..333 get pushed down by one row
# (297->298, 298->299, ... , 333->334), growing the sheet from A1:R333 to A1:R334.
#
# Copy whole rows (not just individual cells) from the bottom up so that number
# formatting (the date style on column D) travels along with the values, and so
# nothing gets overwritten before it has been copied onward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 297
$lastRow  = 333

for ($r = $lastRow; $r -ge $firstRow; $r--) {
    $src = $ws.Range("A" + $r + ":R" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":R" + ($r + 1))
    $src.Copy($dst)
}

# Now overwrite row 297 with the brand-new record. Columns A, B, C, E, F, G, H, I, R
# are identical for every row in this block already (market/category metadata), so
# only the observation-specific columns need to be set.
$ws.Cells.Item($firstRow, 4).Value2  = 45124   # D - Fecha
$ws.Cells.Item($firstRow, 10).Value2 = 200     # J - Volumen
$ws.Cells.Item($firstRow, 11).Value2 = 6000    # K - Precio minimo
$ws.Cells.Item($firstRow, 12).Value2 = 7000    # L - Precio maximo
$ws.Cells.Item($firstRow, 13).Value2 = 6250    # M - Precio promedio ponderado
$ws.Cells.Item($firstRow, 16).Value2 = 125     # P - Precio $/Kg
# N (Unidad de comercializacion), O (Origen) and Q (Kg o Unidades) keep the same
# values the row already had, so they are left untouched.
